# dados_vendas_2023.xlsx — sort the December-2023 ("Dec-23") block of the
# sales table by revenue (faturamento, column B) ascending, then add a
# grand-total row right below the data with a SUM formula.
#
# The sheet has a single table in A1:D901 (header in row 1) grouped into
# monthly blocks (column D = mes-ano). The last block is Dec-23, rows
# 827:901. This commit re-sorts just that block by column B and appends
# row 902 with =SUM(B2:B901) in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the Dec-23 block (A827:D901) by revenue (column B), ascending —
# matches Data > Sort on the filtered/selected block.
$sortRange = $ws.Range("A827:D901")
$sortKey   = $ws.Range("B827:B901")
$sortRange.Sort($sortKey, 1)

# Add the grand-total row right after the last data row.
$ws.Range("B902").Formula = "=SUM(B2:B901)"

# Leave the total cell selected, matching the saved selection state.
$ws.Range("B902").Select()
